$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3854527771472931
$ws.Range("B1").Value = 1.255445003509521
$ws.Range("C1").Value = 6.47393274307251
$ws.Range("D1").Value = 1.739359498023987
$ws.Range("E1").Value = 1.64878237247467
